# fix(publipostage): Try to solve Excel emoji problem
#
# The "statut" column used two emoji glyphs as status markers that were
# causing problems when mail-merged. Replace them with plain symbols:
#   📘 (blue book)  -> ⚠️ (warning sign)
#   📗 (green book) -> ✅ (check mark)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count
$colCount = $usedRange.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        # Keep the string literal on the LEFT of -eq: PowerShell coerces
        # the right-hand side to the left-hand side's type, so a boolean
        # cell value would otherwise be (wrongly) coerced from the emoji
        # string instead of the other way around.
        if ("📘" -eq $val) {
            $cell.Value2 = "⚠️"
        } elseif ("📗" -eq $val) {
            $cell.Value2 = "✅"
        }
    }
}
